$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-10 (columns A-T), reflecting the updated TPM-based
# NATMI ligand-receptor pair statistics for Tnc-Ptprz1 across all sending/target
# cluster combinations (ECs, FAPs, MuSCs).
$data = @(
    @("ECs","Tnc","Ptprz1","ECs",3,1,0.1347866666666667,0.40436,0.03419045085634245,0.03419045085634244,2,0.6666666666666666,0.06163433333333333,0.184903,0.006690894379667537,0.006690894379667537,0.00830748634222222,0.07476737707999999,0.0002287646954730008,0.0002287646954730008),
    @("ECs","Tnc","Ptprz1","FAPs",3,1,0.1347866666666667,0.40436,0.03419045085634245,0.03419045085634244,3,1,0.1030763333333333,0.309229,0.01118975126488057,0.01118975126488057,0.01389331538222222,0.12503983844,0.000382582640716595,0.0003825826407165949),
    @("ECs","Tnc","Ptprz1","MuSCs",3,1,0.1347866666666667,0.40436,0.03419045085634245,0.03419045085634244,3,1,9.046962666666667,27.140888,0.9821193543554519,0.9821193543554518,1.219409941297778,10.97468947168,0.03357910352015285,0.03357910352015284),
    @("FAPs","Tnc","Ptprz1","ECs",3,1,1.308268,3.924804,0.3318597741685039,0.3318597741685039,2,0.6666666666666666,0.06163433333333333,0.184903,0.006690894379667537,0.006690894379667537,0.08063422600133333,0.7257080340119999,0.002220438697821781,0.002220438697821781),
    @("FAPs","Tnc","Ptprz1","FAPs",3,1,1.308268,3.924804,0.3318597741685039,0.3318597741685039,3,1,0.1030763333333333,0.309229,0.01118975126488057,0.01118975126488057,0.1348514684573333,1.213663216116,0.003713428327764998,0.003713428327764998),
    @("FAPs","Tnc","Ptprz1","MuSCs",3,1,1.308268,3.924804,0.3318597741685039,0.3318597741685039,3,1,9.046962666666667,27.140888,0.9821193543554519,0.9821193543554518,11.83585175399467,106.522665785952,0.3259259071429172,0.3259259071429171),
    @("MuSCs","Tnc","Ptprz1","ECs",3,1,2.499176666666667,7.49753,0.6339497749751537,0.6339497749751537,2,0.6666666666666666,0.06163433333333333,0.184903,0.006690894379667537,0.006690894379667537,0.1540350877322222,1.38631578959,0.004241690986372755,0.004241690986372755),
    @("MuSCs","Tnc","Ptprz1","FAPs",3,1,2.499176666666667,7.49753,0.6339497749751537,0.6339497749751537,3,1,0.1030763333333333,0.309229,0.01118975126488057,0.01118975126488057,0.2576059671522222,2.31845370437,0.007093740296398981,0.007093740296398981),
    @("MuSCs","Tnc","Ptprz1","MuSCs",3,1,2.499176666666667,7.49753,0.6339497749751537,0.6339497749751537,3,1,9.046962666666667,27.140888,0.9821193543554519,0.9821193543554518,22.60995800073778,203.48962200664,0.622614343692382,0.6226143436923819)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $rowVals[$c]
    }
}
